$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old contents from A2 and A3
$ws.Range("A2").ClearContents()
$ws.Range("A3").ClearContents()

# Set the new content in A1
$ws.Range("A1").Value = "Hello World !"
